$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$labelText = "Meta description"
$restText = ": Get a taste of classic slot gaming with modern features. Play Cash Spin for free and activate its exciting bonuses with cash prizes and free spins."

$insPoint = $d.Range($metaStart, $metaStart)
$insPoint.Text = ($labelText + $restText)

# Make just the "Meta description" label bold.
$labelRange = $d.Range($metaStart, $metaStart + $labelText.Length)
$labelRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document ("Play Cash Spin Slot Game for Free - ...").
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($i -ne 1 -and $p.Range.Text -like "Play Cash Spin Slot Game for Free - Exciting Bonuses and Classic Graphics*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ------------------------------------------------------------------
# 3) Replace the italic meta-description-like text at the very end of
#    the document with the new image-prompt text.
# ------------------------------------------------------------------
$tailPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Get a taste of classic slot gaming with modern features*") {
        $tailPara = $p
    }
}

$newTail = 'Create a feature image for the game "Cash Spin" that fits the following criteria: - In cartoon style - Features a happy Maya warrior with glasses The image should be lively and colorful, with a cartoon-style depiction of a Maya warrior wearing a big smile and black-rimmed glasses. The warrior should be holding a wheel of fortune in one hand, indicating the game''s bonus feature, while holding a money bag in the other hand to represent the Coin Purse function. The background of the image can be a colorful mixture of the game''s classic symbols, such as diamonds, rubies, emeralds, and dollar signs, arranged in a fun and playful way. Overall, the image should be eye-catching and playful, reflecting the game''s fun and nostalgic atmosphere.'

if ($tailPara -ne $null) {
    $tailStart = $tailPara.Range.Start
    $tailEnd = $tailPara.Range.End - 1   # exclude the paragraph mark
    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Text = $newTail
}
